# BTC_neg_50.xlsx — "added final result calcs after adjusting ML labels"
#
# The sheet "Sheet2" (the ML-labels / back-test calc sheet) has a column C
# of 0/1 "Label" inputs that drive a chain of D (shares bought) / E (running
# portfolio value) formulas down to row 57, which in turn feed the summary
# cells in H2:J3. This commit just flips 20 of those Label cells (simulating
# the ML model's updated predictions) and lets everything downstream
# recalculate. It also fixes/extends two formulas (E50 and E51) whose shared
# range had been mis-split at row 51 with stale absolute row references
# (C57/C58 instead of C51/C52), and moves the active selection to J3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Sheet2" — already the tabSelected sheet

# --- Flip the Label (column C) inputs for the rows that changed ----------
$labelChanges = @{
    4  = 1
    7  = 0
    9  = 1
    10 = 0
    13 = 0
    24 = 1
    25 = 0
    28 = 1
    29 = 1
    31 = 1
    33 = 0
    34 = 1
    37 = 1
    40 = 0
    41 = 0
    42 = 1
    43 = 1
    51 = 1
    52 = 1
    56 = 0
}

foreach ($row in $labelChanges.Keys) {
    $ws.Range("C$row").Value = $labelChanges[$row]
}

# --- Fix the E50/E51 formulas (and propagate the same pattern down to    --
# --- E57) so the "E{n+1} label" reference lines up with the row below    --
# --- instead of the stale out-of-range row references left over from a  --
# --- earlier edit. D51:D57 are rewritten with the equivalent/explicit   --
# --- formula too so the whole D4:D57 / E4:E57 run follows one pattern.  --
$ws.Range("E50").Formula = "=IF(C51=0,IF(D50=0,E49,D50*B50),E49)"

$ws.Range("D51").Formula = "=IF(C51=1,IF(D50=0,E50/B50,D50),0)"
$ws.Range("E51").Formula = "=IF(C52=0,IF(D51=0,E50,D51*B51),E50)"

$ws.Range("D52").Formula = "=IF(C52=1,IF(D51=0,E51/B51,D51),0)"
$ws.Range("E52").Formula = "=IF(C53=0,IF(D52=0,E51,D52*B52),E51)"

$ws.Range("D53").Formula = "=IF(C53=1,IF(D52=0,E52/B52,D52),0)"
$ws.Range("E53").Formula = "=IF(C54=0,IF(D53=0,E52,D53*B53),E52)"

$ws.Range("D54").Formula = "=IF(C54=1,IF(D53=0,E53/B53,D53),0)"
$ws.Range("E54").Formula = "=IF(C55=0,IF(D54=0,E53,D54*B54),E53)"

$ws.Range("D55").Formula = "=IF(C55=1,IF(D54=0,E54/B54,D54),0)"
$ws.Range("E55").Formula = "=IF(C56=0,IF(D55=0,E54,D55*B55),E54)"

$ws.Range("D56").Formula = "=IF(C56=1,IF(D55=0,E55/B55,D55),0)"
$ws.Range("E56").Formula = "=IF(C57=0,IF(D56=0,E55,D56*B56),E55)"

$ws.Range("D57").Formula = "=IF(C57=1,IF(D56=0,E56/B56,D56),0)"
$ws.Range("E57").Formula = "=IF(C58=0,IF(D57=0,E56,D57*B57),E56)"

# --- Move the active selection to J3, as in the saved file ---------------
$ws.Range("J3").Select()
